$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update latitude/longitude for Abingdon (row 2)
$ws.Range("C2").Value = 51.67167
$ws.Range("D2").Value = -1.27833

# Update latitude/longitude for Bagley Wood (row 3)
$ws.Range("C3").Value = 51.718
$ws.Range("D3").Value = -1.2611

# Update longitude for Radcliffe (row 12)
$ws.Range("D12").Value = -2.3277
